# "created tables and plots to analysis page"
# Extends the model-comparison table on Sheet1 with a second day (2019-01-31,
# serial 43496) of results for the same seven cities, re-colors the
# "best model" highlighting that now spans the larger table, and bumps the
# original block's font size to match the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Existing rows 2-8: bump the data font from 10pt -> 12pt (Arial Unicode MS)
#    across the whole C:E block in one shot so existing fill colors (the
#    "1st/2nd/3rd place" highlighting) are preserved on every cell.
# ---------------------------------------------------------------------------
$ws.Range("C2:E8").Font.Size = 12
$ws.Range("C2:E8").Font.Name = "Arial Unicode MS"

# Row 8's "E" cell (Salvador / DTM final) moves from the red (3rd place) to
# the yellow (2nd place) highlight now that a second day of data exists.
$ws.Range("E8").Interior.Color = 65535

# ---------------------------------------------------------------------------
# 2) Append the second day of results as rows 9-15 (same city order as
#    rows 2-8: Manly, Nice, Kauai, Kyoto, Irvine, Amsterdam, Salvador).
# ---------------------------------------------------------------------------
$cities = @("Manly", "Nice", "Kauai", "Kyoto", "Irvine", "Amsterdam", "Salvador")
$newRow = 9
foreach ($city in $cities) {
    $ws.Range("A$newRow").Value = 43496
    $ws.Range("B$newRow").Value = $city
    $ws.Rows.Item($newRow).RowHeight = 18
    $newRow = $newRow + 1
}

# Column A/B formatting matches the block above (date format / plain text).
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A9:A15").PasteSpecial(-4122) | Out-Null
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B9:B15").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Ridge Final / Linear Final / DTM final values for 2019-01-31.
$ws.Range("C9").Value = 0.54200000000000004
$ws.Range("D9").Value = 0.41099999999999998
$ws.Range("E9").Value = 0.35699999999999998

$ws.Range("C10").Value = -0.16600000000000001
$ws.Range("D10").Value = -0.193
$ws.Range("E10").Value = -0.3

$ws.Range("C11").Value = 0.128
$ws.Range("D11").Value = 0.36699999999999999
$ws.Range("E11").Value = 0.052

$ws.Range("C12").Value = -0.90200000000000002
$ws.Range("D12").Value = -1.3620000000000001
$ws.Range("E12").Value = -1.899

$ws.Range("C13").Value = -1.758
$ws.Range("D13").Value = -0.63500000000000001
$ws.Range("E13").Value = -0.067

$ws.Range("C14").Value = 0.20799999999999999
$ws.Range("D14").Value = 0.36899999999999999
$ws.Range("E14").Value = 0.27600000000000002

$ws.Range("C15").Value = -1.3879999999999999
$ws.Range("D15").Value = -1.069
$ws.Range("E15").Value = -1.8380000000000001

# Small (10pt) Arial Unicode MS font for the whole new block, matching the
# font the first block used before step 1 bumped it to 12pt.
$ws.Range("C9:E15").Font.Size = 10
$ws.Range("C9:E15").Font.Name = "Arial Unicode MS"

# ---------------------------------------------------------------------------
# 3) "Best model" highlight colors for the new rows (blue = 1st, yellow =
#    2nd, red never used here - 3rd place is left unfilled in this block).
# ---------------------------------------------------------------------------
$ws.Range("C9").Interior.ThemeColor = 8
$ws.Range("D9").Interior.Color = 65535

$ws.Range("C10").Interior.ThemeColor = 8
$ws.Range("D10").Interior.Color = 65535

$ws.Range("C11").Interior.Color = 65535
$ws.Range("D11").Interior.ThemeColor = 8

$ws.Range("C12").Interior.ThemeColor = 8
$ws.Range("D12").Interior.Color = 65535

$ws.Range("D13").Interior.ThemeColor = 8
$ws.Range("E13").Interior.Color = 65535

$ws.Range("D14").Interior.ThemeColor = 8
$ws.Range("E14").Interior.Color = 65535

$ws.Range("C15").Interior.Color = 65535
$ws.Range("D15").Interior.ThemeColor = 8

# Leave the cursor where the author's last save left it.
$ws.Range("E15").Select() | Out-Null
